$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance row (row 7), mirroring the formatting of the existing rows:
#  - B7 gets the "name/time slot" style used by B4:B6 (Good style w/ red fill override)
#  - C7:H7 get the plain "Good" style used across the other attendance columns
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

$ws.Range("C3:H3").Copy()
$ws.Range("C7:H7").PasteSpecial(-4122)

$ws.Range("B7").Value = "15:30 - 16:05"

$excel.CutCopyMode = $false

# Match the author's final selection (cell B7 was last clicked/edited)
$ws.Range("B7").Select()
